$d = $word.ActiveDocument

# --- 1. "Making data findable" paragraph: metadata of original dataset ---
$d.Content.Find.Execute(
    "Metadata of the original dataset is provided by the publisher and can be found in the repository as well. Machine-readable metadata of both datasets will be created and improve the findability by providing important keywords.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Metadata of the original dataset is provided by the publisher and can be found in the repository as well as machine-readable metadata. For the preprocessed dataset, machine-readable and human-readable metadata files will be created and will improve the findability by providing important keywords.",
    2)

# --- 2. Additionally, the repository will be synced to Zenodo... (append sentence) ---
$d.Content.Find.Execute(
    "Additionally, the repository will be synced to a research data repository (Zenodo) to further improve findability.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Additionally, the repository will be synced to a research data repository (Zenodo) to further improve findability. Zenodo also provides a unique and persistent DOI for the repository.",
    2)

# --- 3. GitHub.com accessibility paragraph (append two sentences) ---
$d.Content.Find.Execute(
    "GitHub.com is a website used by millions of users and therefore a trusted destination for people looking for data. The repository there will be publicly available for everyone and can be downloaded using ssh, https or just with a regular browser. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "GitHub.com is a website used by millions of users and therefore a trusted destination for people looking for data. The repository there will be publicly available for everyone and can be downloaded using ssh, https or just with a regular browser. These methods are free and universally implementable. The same goes for the synced repository on Zenodo.",
    2)

# --- 4. "Zenodo creates a persistent identifier." -> replaced with file naming sentence ---
$d.Content.Find.Execute(
    "Zenodo creates a persistent identifier.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Also, a recognisable file name will be chosen that allows people to know the general contents of the data file without having to open it.",
    2)

# --- 5. All metadata human/machine readable formats (append dublincore sentence) ---
$d.Content.Find.Execute(
    "All metadata will be provided in human-readable (md/pdf) and machine-readable (xml) formats.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "All metadata will be provided in human-readable (md/pdf) and machine-readable (xml) formats and the machine-readable metadata will follow the standard provided by the “dublincore-generator”.",
    2)

# --- 6. Provenance information paragraph (append Metadata sentence) ---
$d.Content.Find.Execute(
    "Provenance information about the data will be provided in a markdown file and the quality of the data is examined in a separate jupyter notebook.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Provenance information about the data will be provided in a markdown file and the quality of the data is examined in a separate jupyter notebook. Metadata will describe the dataset in high level of detail.",
    2)

# --- 7. Remove the now-empty paragraph that followed the Provenance paragraph ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -match "Metadata will describe the dataset in high level of detail\.") {
        $next = $d.Paragraphs.Item($i + 1)
        if ($next.Range.Text.Length -le 1) {
            $next.Range.Delete()
        }
        break
    }
}
